$wb = $excel.ActiveWorkbook
$wsLeft = $wb.Worksheets.Item("left")
$wsRight = $wb.Worksheets.Item("right")

$wsRight.Activate()
$wsRight.Range("A1").Value = "display_id"
$wsRight.Range("B1").Value = "right"
$wsRight.Range("A1:B1").Select()

$wsLeft.Activate()
$wsLeft.Range("A1").Value = "display_id"
$wsLeft.Range("B1").Value = "left"
$wsLeft.Range("B2").Select()
